$d = $word.ActiveDocument

# Locate the existing bullet paragraph that the new bullet must be
# inserted immediately before: "create auth tables in DB;". That
# paragraph's pPr/rPr (auto color, SimSun/Mangal east-asian/complex
# fonts, sz 24, en-US/zh-CN/hi-IN languages) is exactly what the new
# paragraph needs to inherit, so inserting right before it (rather than
# after the previous, green-colored "...CSV files;" bullet) picks up the
# correct formatting automatically.
$anchorText = "create auth tables in DB;"
$rng = $d.Content
$found = $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find anchor paragraph '$anchorText'"
}

$anchorIndex = $rng.Paragraphs.Item(1).Index
$anchorPara = $d.Paragraphs.Item($anchorIndex)

# Insert a new empty paragraph right before it; the new paragraph picks
# up the anchor paragraph's formatting.
$anchorPara.Range.InsertParagraphBefore()

# The freshly created paragraph now occupies the anchor's old slot.
$newPara = $d.Paragraphs.Item($anchorIndex)
$newPara.Range.Text = "localization, default is DE(,) or US(.);"
